$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "2025-01-20"
$ws.Range("B34").Value = "09:07:41"
$ws.Range("C34").Value = "Monday"
$ws.Range("D34").Value = "03"
$ws.Range("E34").Value = 126522
$ws.Range("F34").Value = 142185
$ws.Range("G34").Value = 168971
$ws.Range("H34").Value = 158457
$ws.Range("I34").Value = -1
$ws.Range("J34").Value = 142979
$ws.Range("K34").Value = -1
$ws.Range("L34").Value = -1
$ws.Range("M34").Value = 192500
$ws.Range("N34").Value = 115714
$ws.Range("O34").Value = 45498
$ws.Range("P34").Value = 28506
$ws.Range("Q34").Value = 65765
$ws.Range("R34").Value = -1
$ws.Range("S34").Value = 48715
$ws.Range("T34").Value = -1
